$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new row 43 by copying row 2's formatting/static columns (identical across all rows)
$ws.Range("A2:R2").Copy($ws.Range("A43:R43"))

# Apply the shuffled Fecha/Volumen/Precio values for rows 2-43
$ws.Range("D2").Value2 = 44424
$ws.Range("J2").Value2 = 300
$ws.Range("K2").Value2 = 2500
$ws.Range("L2").Value2 = 3000
$ws.Range("M2").Value2 = 2750
$ws.Range("P2").Value2 = 1375

$ws.Range("D3").Value2 = 44305
$ws.Range("J3").Value2 = 300
$ws.Range("K3").Value2 = 900
$ws.Range("L3").Value2 = 1000
$ws.Range("M3").Value2 = 950
$ws.Range("P3").Value2 = 475

$ws.Range("D4").Value2 = 44243
$ws.Range("J4").Value2 = 200
$ws.Range("K4").Value2 = 2900
$ws.Range("L4").Value2 = 3000
$ws.Range("M4").Value2 = 2950
$ws.Range("P4").Value2 = 1475

$ws.Range("D5").Value2 = 44356
$ws.Range("J5").Value2 = 300
$ws.Range("K5").Value2 = 2400
$ws.Range("L5").Value2 = 2500
$ws.Range("M5").Value2 = 2450
$ws.Range("P5").Value2 = 1225

$ws.Range("D6").Value2 = 44326
$ws.Range("J6").Value2 = 300
$ws.Range("K6").Value2 = 1400
$ws.Range("L6").Value2 = 1500
$ws.Range("M6").Value2 = 1450
$ws.Range("P6").Value2 = 725

$ws.Range("D7").Value2 = 44349
$ws.Range("J7").Value2 = 300
$ws.Range("K7").Value2 = 1800
$ws.Range("L7").Value2 = 2000
$ws.Range("M7").Value2 = 1900
$ws.Range("P7").Value2 = 950

$ws.Range("D8").Value2 = 44298
$ws.Range("J8").Value2 = 300
$ws.Range("K8").Value2 = 1400
$ws.Range("L8").Value2 = 1500
$ws.Range("M8").Value2 = 1450
$ws.Range("P8").Value2 = 725

$ws.Range("D9").Value2 = 44386
$ws.Range("J9").Value2 = 250
$ws.Range("K9").Value2 = 3500
$ws.Range("L9").Value2 = 4000
$ws.Range("M9").Value2 = 3750
$ws.Range("P9").Value2 = 1875

$ws.Range("D10").Value2 = 44266
$ws.Range("J10").Value2 = 300
$ws.Range("K10").Value2 = 1800
$ws.Range("L10").Value2 = 2000
$ws.Range("M10").Value2 = 1900
$ws.Range("P10").Value2 = 950

$ws.Range("D11").Value2 = 44323
$ws.Range("J11").Value2 = 200
$ws.Range("K11").Value2 = 2400
$ws.Range("L11").Value2 = 2500
$ws.Range("M11").Value2 = 2450
$ws.Range("P11").Value2 = 1225

$ws.Range("D12").Value2 = 44165
$ws.Range("J12").Value2 = 300
$ws.Range("K12").Value2 = 1000
$ws.Range("L12").Value2 = 1200
$ws.Range("M12").Value2 = 1100
$ws.Range("P12").Value2 = 550

$ws.Range("D13").Value2 = 44169
$ws.Range("J13").Value2 = 300
$ws.Range("K13").Value2 = 2000
$ws.Range("L13").Value2 = 2500
$ws.Range("M13").Value2 = 2250
$ws.Range("P13").Value2 = 1125

$ws.Range("D14").Value2 = 44405
$ws.Range("J14").Value2 = 300
$ws.Range("K14").Value2 = 3800
$ws.Range("L14").Value2 = 4000
$ws.Range("M14").Value2 = 3900
$ws.Range("P14").Value2 = 1950

$ws.Range("D15").Value2 = 44431
$ws.Range("J15").Value2 = 300
$ws.Range("K15").Value2 = 1900
$ws.Range("L15").Value2 = 2000
$ws.Range("M15").Value2 = 1950
$ws.Range("P15").Value2 = 975

$ws.Range("D16").Value2 = 44181
$ws.Range("J16").Value2 = 250
$ws.Range("K16").Value2 = 1400
$ws.Range("L16").Value2 = 1500
$ws.Range("M16").Value2 = 1450
$ws.Range("P16").Value2 = 725

$ws.Range("D17").Value2 = 44203
$ws.Range("J17").Value2 = 300
$ws.Range("K17").Value2 = 2000
$ws.Range("L17").Value2 = 2500
$ws.Range("M17").Value2 = 2250
$ws.Range("P17").Value2 = 1125

$ws.Range("D18").Value2 = 44312
$ws.Range("J18").Value2 = 300
$ws.Range("K18").Value2 = 1000
$ws.Range("L18").Value2 = 1200
$ws.Range("M18").Value2 = 1100
$ws.Range("P18").Value2 = 550

$ws.Range("D19").Value2 = 44258
$ws.Range("J19").Value2 = 150
$ws.Range("K19").Value2 = 2400
$ws.Range("L19").Value2 = 2500
$ws.Range("M19").Value2 = 2450
$ws.Range("P19").Value2 = 1225

$ws.Range("D20").Value2 = 44410
$ws.Range("J20").Value2 = 250
$ws.Range("K20").Value2 = 2800
$ws.Range("L20").Value2 = 3000
$ws.Range("M20").Value2 = 2900
$ws.Range("P20").Value2 = 1450

$ws.Range("D21").Value2 = 44284
$ws.Range("J21").Value2 = 300
$ws.Range("K21").Value2 = 1800
$ws.Range("L21").Value2 = 2000
$ws.Range("M21").Value2 = 1900
$ws.Range("P21").Value2 = 950

$ws.Range("D22").Value2 = 44343
$ws.Range("J22").Value2 = 300
$ws.Range("K22").Value2 = 1500
$ws.Range("L22").Value2 = 2000
$ws.Range("M22").Value2 = 1750
$ws.Range("P22").Value2 = 875

$ws.Range("D23").Value2 = 44293
$ws.Range("J23").Value2 = 250
$ws.Range("K23").Value2 = 1500
$ws.Range("L23").Value2 = 1800
$ws.Range("M23").Value2 = 1650
$ws.Range("P23").Value2 = 825

$ws.Range("D24").Value2 = 44320
$ws.Range("J24").Value2 = 250
$ws.Range("K24").Value2 = 1400
$ws.Range("L24").Value2 = 1500
$ws.Range("M24").Value2 = 1450
$ws.Range("P24").Value2 = 725

$ws.Range("D25").Value2 = 44221
$ws.Range("J25").Value2 = 200
$ws.Range("K25").Value2 = 2900
$ws.Range("L25").Value2 = 3000
$ws.Range("M25").Value2 = 2950
$ws.Range("P25").Value2 = 1475

$ws.Range("D26").Value2 = 44334
$ws.Range("J26").Value2 = 200
$ws.Range("K26").Value2 = 2800
$ws.Range("L26").Value2 = 3000
$ws.Range("M26").Value2 = 2900
$ws.Range("P26").Value2 = 1450

$ws.Range("D27").Value2 = 44370
$ws.Range("J27").Value2 = 400
$ws.Range("K27").Value2 = 3400
$ws.Range("L27").Value2 = 3500
$ws.Range("M27").Value2 = 3445
$ws.Range("P27").Value2 = 1722

$ws.Range("D28").Value2 = 44267
$ws.Range("J28").Value2 = 300
$ws.Range("K28").Value2 = 1400
$ws.Range("L28").Value2 = 1500
$ws.Range("M28").Value2 = 1450
$ws.Range("P28").Value2 = 725

$ws.Range("D29").Value2 = 44249
$ws.Range("J29").Value2 = 300
$ws.Range("K29").Value2 = 2400
$ws.Range("L29").Value2 = 2500
$ws.Range("M29").Value2 = 2450
$ws.Range("P29").Value2 = 1225

$ws.Range("D30").Value2 = 44272
$ws.Range("J30").Value2 = 250
$ws.Range("K30").Value2 = 2800
$ws.Range("L30").Value2 = 3000
$ws.Range("M30").Value2 = 2900
$ws.Range("P30").Value2 = 1450

$ws.Range("D31").Value2 = 44235
$ws.Range("J31").Value2 = 250
$ws.Range("K31").Value2 = 4500
$ws.Range("L31").Value2 = 5000
$ws.Range("M31").Value2 = 4750
$ws.Range("P31").Value2 = 2375

$ws.Range("D32").Value2 = 44433
$ws.Range("J32").Value2 = 200
$ws.Range("K32").Value2 = 1800
$ws.Range("L32").Value2 = 2000
$ws.Range("M32").Value2 = 1900
$ws.Range("P32").Value2 = 950

$ws.Range("D33").Value2 = 44397
$ws.Range("J33").Value2 = 300
$ws.Range("K33").Value2 = 3500
$ws.Range("L33").Value2 = 4000
$ws.Range("M33").Value2 = 3750
$ws.Range("P33").Value2 = 1875

$ws.Range("D34").Value2 = 44263
$ws.Range("J34").Value2 = 270
$ws.Range("K34").Value2 = 1900
$ws.Range("L34").Value2 = 2000
$ws.Range("M34").Value2 = 1950
$ws.Range("P34").Value2 = 975

$ws.Range("D35").Value2 = 44253
$ws.Range("J35").Value2 = 300
$ws.Range("K35").Value2 = 2400
$ws.Range("L35").Value2 = 2500
$ws.Range("M35").Value2 = 2450
$ws.Range("P35").Value2 = 1225

$ws.Range("D36").Value2 = 44176
$ws.Range("J36").Value2 = 300
$ws.Range("K36").Value2 = 1900
$ws.Range("L36").Value2 = 2000
$ws.Range("M36").Value2 = 1950
$ws.Range("P36").Value2 = 975

$ws.Range("D37").Value2 = 44417
$ws.Range("J37").Value2 = 300
$ws.Range("K37").Value2 = 3000
$ws.Range("L37").Value2 = 3500
$ws.Range("M37").Value2 = 3250
$ws.Range("P37").Value2 = 1625

$ws.Range("D38").Value2 = 44237
$ws.Range("J38").Value2 = 200
$ws.Range("K38").Value2 = 2500
$ws.Range("L38").Value2 = 3000
$ws.Range("M38").Value2 = 2750
$ws.Range("P38").Value2 = 1375

$ws.Range("D39").Value2 = 44428
$ws.Range("J39").Value2 = 270
$ws.Range("K39").Value2 = 3500
$ws.Range("L39").Value2 = 3800
$ws.Range("M39").Value2 = 3650
$ws.Range("P39").Value2 = 1825

$ws.Range("D40").Value2 = 44442
$ws.Range("J40").Value2 = 200
$ws.Range("K40").Value2 = 2400
$ws.Range("L40").Value2 = 2500
$ws.Range("M40").Value2 = 2450
$ws.Range("P40").Value2 = 1225

$ws.Range("D41").Value2 = 44435
$ws.Range("J41").Value2 = 500
$ws.Range("K41").Value2 = 1800
$ws.Range("L41").Value2 = 2000
$ws.Range("M41").Value2 = 1930
$ws.Range("P41").Value2 = 965

$ws.Range("D42").Value2 = 44319
$ws.Range("J42").Value2 = 300
$ws.Range("K42").Value2 = 1900
$ws.Range("L42").Value2 = 2000
$ws.Range("M42").Value2 = 1950
$ws.Range("P42").Value2 = 975

$ws.Range("D43").Value2 = 44279
$ws.Range("J43").Value2 = 200
$ws.Range("K43").Value2 = 1700
$ws.Range("L43").Value2 = 1800
$ws.Range("M43").Value2 = 1750
$ws.Range("P43").Value2 = 875
